$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65

# Row 4 updates
$ws.Range("G4").Value = 2.75
$ws.Range("I4").Value = 2.5
$ws.Range("J4").Value = 3.2
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 3.05
$ws.Range("O4").Value = 1.29
$ws.Range("W4").Value = 9.25
$ws.Range("X4").Value = 15
$ws.Range("Y4").Value = 9.75
$ws.Range("AA4").Value = 22
$ws.Range("AB4").Value = 29
$ws.Range("AD4").Value = 6
$ws.Range("AH4").Value = 8.25
$ws.Range("AI4").Value = 12.5
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 29
$ws.Range("AO4").Value = 14
$ws.Range("AP4").Value = 19
$ws.Range("AQ4").Value = 60
$ws.Range("AR4").Value = 80
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 2.65
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 19
$ws.Range("AZ4").Value = 55
$ws.Range("BA4").Value = 80
$ws.Range("BB4").Value = 200
